# Adds the missing utm_* / external_id header rows to the "sw_signups"
# (signups_table) sheet, turning the previously-blank placeholder rows
# 10-12 into real data rows and appending two brand-new rows (13, 14).
#
# Columns: A=headers  B=null_allowed  C=type  D=minimum  E=maximum
#          F=minLength G=maxLength H=pattern I=enum J=description
#          K=examples

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)          # sw_summary  (used only as a format donor)
$ws  = $wb.Worksheets.Item(2)          # sw_signups  (signups_table)

# ---------------------------------------------------------------------
# 1) Normalise formatting for the two existing rows (8 & 9) so every
#    populated cell shares the plain "data row" look used by rows 5/6,
#    then stamp the new text for row 8 / the K cell of row 9.
# ---------------------------------------------------------------------
$ws.Range("A5:K5").Copy()
$ws.Range("A8:K8").PasteSpecial(-4122)

$ws.Range("K7").Copy()                 # right-aligned "examples" style
$ws.Range("K9").PasteSpecial(-4122)

$ws.Range("A8").Value = "external_id"
$ws.Range("B8").Value = "Yes"
$ws.Range("C8").Value = "string"
$ws.Range("K8").Value = "#ioy2fcf"

$ws.Range("A9").Value = "utm_id"
$ws.Range("B9").Value = "Yes"
$ws.Range("C9").Value = "string"
$ws.Range("K9").Value = "wvninnewn"

# ---------------------------------------------------------------------
# 2) Turn the blank placeholder rows 10-12 into real data rows, and add
#    two brand-new rows 13-14 the same way.
# ---------------------------------------------------------------------
$rows = @(
    @{ Row = 10; A = "utm_source";   K = "youtube" },
    @{ Row = 11; A = "utm_medium";   K = "social" },
    @{ Row = 12; A = "utm_campaign"; K = "fireship" },
    @{ Row = 13; A = "utm_term";     K = "apple" },
    @{ Row = 14; A = "utm_content";  K = "logolink" }
)

foreach ($item in $rows) {
    $r = $item.Row

    # Base formatting for the row = the plain "data row" style (A:K).
    $ws.Range("A5:K5").Copy()
    $ws.Range("A" + $r + ":K" + $r).PasteSpecial(-4122)

    $ws.Range("A" + $r).Value = $item.A
    $ws.Range("B" + $r).Value = "Yes"
    $ws.Range("C" + $r).Value = "string"
    $ws.Range("K" + $r).Value = $item.K
}

# ---------------------------------------------------------------------
# 3) Every data row (8-14) also carries blank, plainly-formatted cells
#    out to column Z (a left-over of how these rows were authored).
# ---------------------------------------------------------------------
$ws1.Range("A9").Copy()
$ws.Range("L8:Z14").PasteSpecial(-4122)

Write-Host "done"
